$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.608.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.105.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.104.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.87%  "

$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("E10").Value = "  +0.29%  "

$ws.Range("E11").Value = "  +2.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.640.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.00%  "

$ws.Range("E16").Value = "  +1.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.689.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.113.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "336.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.513"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.83%  "

$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0923"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.18%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("E32").Value = "  +2.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "

$ws.Range("E34").Value = "  +1.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("E39").Value = "  +2.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0663"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("E41").Value = "  +12.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.148.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.686"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.308.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.69%  "

$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.977"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.07%  "

$ws.Range("E50").Value = "  +0.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.99%  "
